$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("备选")

# Clear the M-column penalty markers on rows 6, 8, 12, 16, 42 (their -1 values are removed)
$ws.Range("M6").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("M42").ClearContents()

# Add a new bond entry (开润转债 / Kairun) as rows 79-80, mirroring the 2-row
# block pattern used by every other entry (odd row = data, even row = score inputs).
$ws.Rows("79:80").Insert()

$ws.Range("A77:P77").Copy()
$ws.Range("A79").PasteSpecial(-4122)
$ws.Range("A78:P78").Copy()
$ws.Range("A80").PasteSpecial(-4122)

# Row 79: main data row for 开润转债 (Kairun Convertible Bond)
$ws.Range("A79").Value = "开润转债"
$ws.Range("C79").Value = "箱包"
$ws.Range("D79").Value = "丁碧霞"
$ws.Range("E79").Value = "小米"
$ws.Range("F79").Value = 1.5
$ws.Range("G79").Value = " 2024-03-31"
$ws.Range("H79").Value = 0.3264
$ws.Range("I79").Value = 59.7197
$ws.Range("J79").Value = 11.45
$ws.Range("K79").Value = "足够还可转债"
$ws.Range("L79").Value = Get-Date -Year 2025 -Month 12 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("M79").Formula = '=DATEDIF(TODAY(),L79,"y")&"年"&DATEDIF(TODAY(),L79,"ym")&"月"'
$ws.Range("P79").Formula = "=SUM(B80:O80)"

# Row 80: score-input row under the new entry
$ws.Range("D80").Value = 1
$ws.Range("E80").Value = 1
$ws.Range("F80").Value = 1
